$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "[Chaomin%Wu%xref no email%0,       Xiaoyan%Chen%xref no email%0,       Yanping%Cai%xref no email%0,       Jia\u2019an%Xia%xref no email%0,       Xing%Zhou%xref no email%0,       Sha%Xu%xref no email%0,       Hanping%Huang%xref no email%0,       Li%Zhang%xref no email%0,       Xia%Zhou%xref no email%0,       Chunling%Du%xref no email%0,       Yuye%Zhang%xref no email%0,       Juan%Song%xref no email%0,       Sijiao%Wang%xref no email%0,       Yencheng%Chao%xref no email%0,       Zeyong%Yang%xref no email%0,       Jie%Xu%xref no email%0,       Xin%Zhou%xref no email%0,       Dechang%Chen%xref no email%0,       Weining%Xiong%xref no email%0,       Lei%Xu%xref no email%0,       Feng%Zhou%xref no email%0,       Jinjun%Jiang%xref no email%0,       Chunxue%Bai%xref no email%0,       Junhua%Zheng%xref no email%0,       Yuanlin%Song%xref no email%0]"
$ws.Range("E3").Value = "[Michael N.%Young%NULL%1,                            Dhaval%Kolte%NULL%1,                            Mary E.%Cadigan%NULL%1,                            Elizabeth%Laikhter%NULL%1,                            Kevin%Sinclair%NULL%1,                            Eugene%Pomerantsev%NULL%1,                            Michael A.%Fifer%NULL%1,                            Thoralf M.%Sundt%NULL%1,                            Robert W.%Yeh%NULL%1,                            Farouc A.%Jaffer%fjaffer@mgh.harvard.edu%1]"
$ws.Range("E4").Value = "[Antonelli%Massimo%coreGivesNoEmail%0,                           Cabrini%Luca%coreGivesNoEmail%1,                           Castelli%Antonio%coreGivesNoEmail%1,                           Cecconi%Maurizio%coreGivesNoEmail%1,                           Cereda%Danilo%coreGivesNoEmail%1,                           Coluccello%Antonio%coreGivesNoEmail%1,                           Foti%Giuseppe%coreGivesNoEmail%1,                           Fumagalli%Roberto%coreGivesNoEmail%1,                           Grasselli%Giacomo%coreGivesNoEmail%1,                           Iotti%Giorgio%coreGivesNoEmail%1,                           Latronico%Nicola%coreGivesNoEmail%1,                           Lorini%Luca%coreGivesNoEmail%1,                           Merler%Stefano%coreGivesNoEmail%1,                           Natalini%Giuseppe%coreGivesNoEmail%1,                           Pesenti%Antonio%coreGivesNoEmail%1,                           Piatti%Alessandra%coreGivesNoEmail%1,                           Ranieri%Marco Vito%coreGivesNoEmail%1,                           Scandroglio%Anna Mara%coreGivesNoEmail%1,                           Storti%Enrico%coreGivesNoEmail%1,                           Zanella%Alberto%coreGivesNoEmail%1,                           Zangrillo%Alberto%coreGivesNoEmail%1]"
$ws.Range("E5").Value = "[Chaolin%Huang%NULL%0,                            Yeming%Wang%NULL%0,                            Xingwang%Li%NULL%0,                            Lili%Ren%NULL%0,                            Jianping%Zhao%NULL%0,                            Yi%Hu%NULL%0,                            Li%Zhang%NULL%0,                            Guohui%Fan%NULL%0,                            Jiuyang%Xu%NULL%0,                            Xiaoying%Gu%NULL%0,                            Zhenshun%Cheng%NULL%0,                            Ting%Yu%NULL%0,                            Jiaan%Xia%NULL%0,                            Yuan%Wei%NULL%0,                            Wenjuan%Wu%NULL%0,                            Xuelei%Xie%NULL%0,                            Wen%Yin%NULL%0,                            Hui%Li%NULL%0,                            Min%Liu%NULL%0,                            Yan%Xiao%NULL%0,                            Hong%Gao%NULL%0,                            Li%Guo%NULL%0,                            Jungang%Xie%NULL%0,                            Guangfa%Wang%NULL%0,                            Rongmeng%Jiang%NULL%0,                            Zhancheng%Gao%NULL%0,                            Qi%Jin%NULL%0,                            Jianwei%Wang%wangjw28@163.com%0,                            Bin%Cao%caobin_ben@163.com%0]"
$ws.Range("E6").Value = "[Suxin%Wan%NULL%0,                            Yi%Xiang%NULL%0,                            Wei%Fang%NULL%0,                            Yu%Zheng%yuzheng1@cdutcm.edu.cn%0,                            Boqun%Li%279685211@qq.com%0,                            Boqun%Li%279685211@qq.com%0,                            Yanjun%Hu%huyanjun@163.com%0,                            Chunhui%Lang%NULL%0,                            Daoqiu%Huang%NULL%0,                            Qiuyan%Sun%NULL%0,                            Yan%Xiong%NULL%0,                            Xia%Huang%NULL%0,                            Jinglong%Lv%NULL%0,                            Yaling%Luo%NULL%0,                            Li%Shen%NULL%0,                            Haoran%Yang%NULL%0,                            Gu%Huang%NULL%0,                            Ruishan%Yang%NULL%0]"
$ws.Range("E7").Value = "[Matt%Arentz%xref no email%0,       Eric%Yim%xref no email%0,       Lindy%Klaff%xref no email%0,       Sharukh%Lokhandwala%xref no email%0,       Francis X.%Riedo%xref no email%0,       Maria%Chong%xref no email%0,       Melissa%Lee%xref no email%0]"
$ws.Range("E8").Value = "[Sara C.%Auld%NULL%0,                            Mark%Caridi-Scheible%NULL%3,                            James M.%Blum%NULL%2,                            Chad%Robichaux%NULL%2,                            Colleen%Kraft%NULL%2,                            Jesse T.%Jacob%NULL%2,                            Craig S.%Jabaley%NULL%2,                            David%Carpenter%NULL%2,                            Roberta%Kaplow%NULL%2,                            Alfonso C.%Hernandez-Romieu%NULL%2,                            Max W.%Adelman%NULL%2,                            Greg S.%Martin%NULL%2,                            Craig M.%Coopersmith%NULL%2,                            David J.%Murphy%NULL%2,                            NULL%NULL%NULL%0,                             S.%Auld%null%1,                             M.% Caridi-Scheible%null%1,                             J. M.% Blum%null%1,                             C. J.% Robichaux%null%1,                             C. S.% Kraft%null%1,                             J. T.% Jacob%null%1,                             C. S.% Jabaley%null%1,                             D.% Carpenter%null%1,                             R.% Kaplow%null%1,                             A. C.% Hernandez%null%1,                             M. W.% Adelman%null%1,                             G. S.% Martin%null%1,                             C. M.% Coopersmith%null%1,                             D. J.% Murphy%null%1,                               % Emory COVID-19 Quality and Clinical Research Collaborative%null%1,                           S.%Auld%null%1,                           M.% Caridi-Scheible%null%1,                           J. M.% Blum%null%1,                           C. J.% Robichaux%null%1,                           C. S.% Kraft%null%1,                           J. T.% Jacob%null%1,                           C. S.% Jabaley%null%1,                           D.% Carpenter%null%1,                           R.% Kaplow%null%1,                           A. C.% Hernandez%null%1,                           M. W.% Adelman%null%1,                           G. S.% Martin%null%1,                           C. M.% Coopersmith%null%1,                           D. J.% Murphy%null%1,                             % Emory COVID-19 Quality and Clinical Research Collaborative%null%1,         S.%Auld%null%1,         M.% Caridi-Scheible%null%1,         J. M.% Blum%null%1,         C. J.% Robichaux%null%1,         C. S.% Kraft%null%1,         J. T.% Jacob%null%1,         C. S.% Jabaley%null%1,         D.% Carpenter%null%1,         R.% Kaplow%null%1,         A. C.% Hernandez%null%1,         M. W.% Adelman%null%1,         G. S.% Martin%null%1,         C. M.% Coopersmith%null%1,         D. J.% Murphy%null%1,           % Emory COVID-19 Quality and Clinical Research Collaborative%null%1]"
$ws.Range("E9").Value = "[Pavan K.%Bhatraju%NULL%0,                            Bijan J.%Ghassemieh%NULL%2,                            Michelle%Nichols%NULL%2,                            Richard%Kim%NULL%2,                            Keith R.%Jerome%NULL%2,                            Arun K.%Nalla%NULL%2,                            Alexander L.%Greninger%NULL%2,                            Sudhakar%Pipavath%NULL%2,                            Mark M.%Wurfel%NULL%2,                            Laura%Evans%NULL%2,                            Patricia A.%Kritek%NULL%2,                            T. Eoin%West%NULL%2,                            Andrew%Luks%NULL%2,                            Anthony%Gerbino%NULL%2,                            Chris R.%Dale%NULL%2,                            Jason D.%Goldman%NULL%2,                            Shane%O’Mahony%NULL%2,                            Carmen%Mikacenic%NULL%2]"
$ws.Range("E10").Value = "[Jun%Chen%NULL%0,                            Tangkai%Qi%NULL%1,                            Li%Liu%NULL%1,                            Yun%Ling%NULL%0,                            Zhiping%Qian%NULL%1,                            Tao%Li%NULL%0,                            Feng%Li%NULL%1,                            Qingnian%Xu%NULL%1,                            Yuyi%Zhang%NULL%1,                            Shuibao%Xu%NULL%1,                            Zhigang%Song%NULL%1,                            Yigang%Zeng%NULL%1,                            Yinzhong%Shen%NULL%1,                            Yuxin%Shi%NULL%2,                            Tongyu%Zhu%NULL%2,                            Hongzhou%Lu%NULL%4]"
$ws.Range("E11").Value = "[Nanshan%Chen%NULL%0,                            Min%Zhou%NULL%0,                            Xuan%Dong%NULL%0,                            Jieming%Qu%NULL%0,                            Fengyun%Gong%NULL%0,                            Yang%Han%NULL%0,                            Yang%Qiu%NULL%0,                            Jingli%Wang%NULL%0,                            Ying%Liu%NULL%0,                            Yuan%Wei%NULL%0,                            Jia'an%Xia%NULL%0,                            Ting%Yu%NULL%0,                            Xinxin%Zhang%NULL%0,                            Li%Zhang%NULL%0]"
$ws.Range("E12").Value = "[Wei-jie%Guan%NULL%0,                            Zheng-yi%Ni%NULL%0,                            Zheng-yi%Ni%NULL%0,                            Yu%Hu%NULL%0,                            Wen-hua%Liang%NULL%0,                            Chun-quan%Ou%NULL%0,                            Jian-xing%He%NULL%0,                            Lei%Liu%NULL%0,                            Hong%Shan%NULL%0,                            Chun-liang%Lei%NULL%0,                            David S.C.%Hui%NULL%0,                            Bin%Du%NULL%0,                            Lan-juan%Li%NULL%0,                            Guang%Zeng%NULL%0,                            Kwok-Yung%Yuen%NULL%0,                            Ru-chong%Chen%NULL%0,                            Chun-li%Tang%NULL%0,                            Tao%Wang%NULL%0,                            Ping-yan%Chen%NULL%0,                            Jie%Xiang%NULL%0,                            Shi-yue%Li%NULL%0,                            Jin-lin%Wang%NULL%0,                            Zi-jing%Liang%NULL%0,                            Yi-xiang%Peng%NULL%0,                            Li%Wei%NULL%0,                            Yong%Liu%NULL%0,                            Ya-hua%Hu%NULL%0,                            Peng%Peng%NULL%0,                            Jian-ming%Wang%NULL%0,                            Ji-yang%Liu%NULL%0,                            Zhong%Chen%NULL%0,                            Gang%Li%NULL%0,                            Zhi-jian%Zheng%NULL%0,                            Shao-qin%Qiu%NULL%0,                            Jie%Luo%NULL%0,                            Chang-jiang%Ye%NULL%0,                            Shao-yong%Zhu%NULL%0,                            Nan-shan%Zhong%NULL%0]"
$ws.Range("E13").Value = "[Pingzheng%Mo%NULL%0,                            Yuanyuan%Xing%NULL%2,                            Yu%Xiao%NULL%2,                            Liping%Deng%NULL%3,                            Qiu%Zhao%NULL%2,                            Hongling%Wang%NULL%2,                            Yong%Xiong%NULL%4,                            Zhenshun%Cheng%NULL%0,                            Shicheng%Gao%NULL%3,                            Ke%Liang%NULL%3,                            Mingqi%Luo%NULL%3,                            Tielong%Chen%NULL%2,                            Shihui%Song%NULL%3,                            Zhiyong%Ma%NULL%3,                            Xiaoping%Chen%NULL%3,                            Ruiying%Zheng%NULL%2,                            Qian%Cao%NULL%2,                            Fan%Wang%fanndywang@foxmail.com%2,                            Yongxi%Zhang%znact1936@126.com%0]"
$ws.Range("E14").Value = "[G -Q%Qian%NULL%0,                            N -B%Yang%NULL%1,                            F%Ding%NULL%1,                            A H Y%Ma%NULL%1,                            Z -Y%Wang%NULL%2,                            Y -F%Shen%NULL%1,                            C -W%Shi%NULL%1,                            X%Lian%NULL%1,                            J -G%Chu%NULL%1,                            L%Chen%chxmin@hotmail.com%1,                            Z -Y%Wang%NULL%0,                            D -W%Ren%NULL%1,                            G -X%Li%NULL%1,                            X -Q%Chen%chxmin@hotmail.com%1,                            H -J%Shen%NULL%1,                            X -M%Chen%chxmin@hotmail.com%1]"
$ws.Range("E15").Value = "[Flavien%Vincent%NULL%1,                            Hugues%Spillemaeker%NULL%1,                            Maéva%Kyheng%NULL%1,                            Cassandre%Belin‐Vincent%NULL%1,                            Cédric%Delhaye%NULL%1,                            Adeline%Piérache%NULL%1,                            Tom%Denimal%NULL%1,                            Basile%Verdier%NULL%1,                            Nicolas%Debry%NULL%1,                            Mouhamed%Moussa%NULL%1,                            Guillaume%Schurtz%NULL%1,                            Sina%Porouchani%NULL%1,                            Alessandro%Cosenza%NULL%1,                            Francis%Juthier%NULL%1,                            Thibault%Pamart%NULL%1,                            Marjorie%Richardson%NULL%1,                            Augustin%Coisne%NULL%1,                            Adrien%Hertault%NULL%1,                            Jonathan%Sobocinski%NULL%1,                            Thomas%Modine%NULL%1,                            François%Pontana%NULL%1,                            Alain%Duhamel%NULL%2,                            Julien%Labreuche%NULL%2,                            Eric%Van Belle%ericvanbelle@aol.com%1]"
$ws.Range("E16").Value = "[Carb\u00f3-Cisnero%Yacquelin%coreGivesNoEmail%0,                           Fern\u00e1ndez-Gonz\u00e1lez%Paula%coreGivesNoEmail%0,                           Hierrezuelo-Rojas%Naifi%coreGivesNoEmail%0,                           Subert-Salas%Lizandra%coreGivesNoEmail%0]"
$ws.Range("E17").Value = "[Dawei%Wang%wdw_syr@163.com%1,                            Yimei%Yin%18256985@qq.com%1,                            Chang%Hu%1391369733@qq.com%3,                            Xing%Liu%375117565@qq.com%1,                            Xingguo%Zhang%zhangxinguo7909@163.com%1,                            Shuliang%Zhou%350925556@qq.com%1,                            Mingzhi%Jian%Jmzhak@aliyun.com%1,                            Haibo%Xu%xuhaibo1120@hotmail.com%0,                            John%Prowle%j.prowle@qmul.ac.uk%1,                            Bo%Hu%hobbier1979@163.com%0,                            Yirong%Li%liyirong838@163.com%1,                            Zhiyong%Peng%Pengzy5@hotmail.com%0]"
$ws.Range("E18").Value = "[Wentao%Xu%xref no email%1,       Shuang%Qu%xref no email%1,       Mengying%Xing%xref no email%1,       Ming%Zhang%xref no email%1,       Geng%Lu%xref no email%1,       Zhicong%Liao%xref no email%1,       Kimberly%Griffin%xref no email%1,       Jun%Wang%xref no email%1,       Ke%Zen%xref no email%1,       Bing%Yao%xref no email%1,       Hongwei%Liang%xref no email%1,       Jianping%Zhang%xref no email%1]"
$ws.Range("E19").Value = "[Xiao-Wei%Xu%NULL%0,                            Xiao-Xin%Wu%NULL%2,                            Xian-Gao%Jiang%NULL%2,                            Kai-Jin%Xu%NULL%2,                            Ling-Jun%Ying%NULL%2,                            Chun-Lian%Ma%NULL%2,                            Shi-Bo%Li%NULL%2,                            Hua-Ying%Wang%NULL%2,                            Sheng%Zhang%NULL%2,                            Hai-Nv%Gao%NULL%2,                            Ji-Fang%Sheng%NULL%2,                            Hong-Liu%Cai%NULL%2,                            Yun-Qing%Qiu%NULL%2,                            Lan-Juan%Li%NULL%2]"
$ws.Range("E20").Value = "[Guqin%Zhang%NULL%0,                            Chang%Hu%NULL%0,                            Linjie%Luo%NULL%2,                            Fang%Fang%NULL%0,                            Yongfeng%Chen%NULL%2,                            Jianguo%Li%NULL%2,                            Zhiyong%Peng%NULL%0,                            Huaqin%Pan%phq2012@whu.edu.cn%2]"
$ws.Range("E21").Value = "[Fang%Zheng%fangzheng99@sina.cn%1,                            Chun%Liao%lc7spring@163.com%1,                            Qi-hong%Fan%1187286305@qq.com%1,                            Hong-bo%Chen%hbchen@hust.edu.cn%1,                            Xue-gong%Zhao%NULL%1,                            Zhong-guo%Xie%NULL%1,                            Xi-lin%Li%NULL%1,                            Chun-xi%Chen%NULL%1,                            Xiao-xia%Lu%NULL%1,                            Zhi-sheng%Liu%NULL%1,                            Wei%Lu%NULL%1,                            Chun-bao%Chen%NULL%1,                            Rong%Jiao%NULL%1,                            Ai-ming%Zhang%NULL%1,                            Jin-tang%Wang%NULL%1,                            Xi-wei%Ding%NULL%1,                            Yao-guang%Zeng%NULL%1,                            Li-ping%Cheng%NULL%1,                            Qing-feng%Huang%NULL%1,                            Jiang%Wu%NULL%0,                            Xi-chang%Luo%NULL%1,                            Zhu-jun%Wang%NULL%1,                            Yan-yan%Zhong%NULL%1,                            Yan%Bai%yanbaixh@hust.edu.cn%1,                            Xiao-yan%Wu%xwu@hust.edu.cn%1,                            Run-ming%Jin%jinrunm@qq.com%1]"
$ws.Range("E22").Value = "[Fei%Zhou%NULL%0,                            Ting%Yu%NULL%0,                            Ronghui%Du%NULL%0,                            Guohui%Fan%NULL%0,                            Ying%Liu%NULL%0,                            Zhibo%Liu%NULL%0,                            Jie%Xiang%NULL%0,                            Yeming%Wang%NULL%0,                            Bin%Song%NULL%0,                            Xiaoying%Gu%NULL%0,                            Lulu%Guan%NULL%0,                            Yuan%Wei%NULL%0,                            Hui%Li%NULL%0,                            Xudong%Wu%NULL%0,                            Jiuyang%Xu%NULL%0,                            Shengjin%Tu%NULL%0,                            Yi%Zhang%NULL%0,                            Hua%Chen%NULL%0,                            Bin%Cao%NULL%0]"
$ws.Range("E23").Value = "[Yihui%Huang%NULL%0,                            Mengqi%Tu%NULL%2,                            Shipei%Wang%NULL%2,                            Sichao%Chen%NULL%2,                            Wei%Zhou%NULL%0,                            Danyang%Chen%NULL%2,                            Lin%Zhou%NULL%2,                            Min%Wang%NULL%2,                            Yan%Zhao%NULL%4,                            Wen%Zeng%NULL%2,                            Qi%Huang%NULL%2,                            Hai'bo%Xu%NULL%2,                            Zeming%Liu%NULL%2,                            Liang%Guo%NULL%2]"
$ws.Range("E24").Value = "[Xiaobo%Yang%NULL%0,                            Yuan%Yu%NULL%0,                            Jiqian%Xu%NULL%0,                            Huaqing%Shu%NULL%0,                            Jia'an%Xia%NULL%0,                            Hong%Liu%NULL%0,                            Yongran%Wu%NULL%0,                            Lu%Zhang%NULL%0,                            Zhui%Yu%NULL%0,                            Minghao%Fang%NULL%0,                            Ting%Yu%NULL%0,                            Yaxin%Wang%NULL%0,                            Shangwen%Pan%NULL%0,                            Xiaojing%Zou%NULL%0,                            Shiying%Yuan%NULL%0,                            You%Shang%NULL%0]"
$ws.Range("E25").Value = "[Wen-hua%Liang%NULL%2,                            Wei-jie%Guan%NULL%2,                            Cai-chen%Li%NULL%3,                            Yi-min%Li%NULL%3,                            Heng-rui%Liang%NULL%3,                            Yi%Zhao%NULL%3,                            Xiao-qing%Liu%NULL%3,                            Ling%Sang%NULL%0,                            Ru-chong%Chen%NULL%0,                            Chun-li%Tang%NULL%0,                            Tao%Wang%NULL%0,                            Wei%Wang%NULL%0,                            Qi-hua%He%NULL%2,                            Zi-sheng%Chen%NULL%3,                            Sook-San%Wong%NULL%2,                            Mark%Zanin%NULL%2,                            Jun%Liu%NULL%2,                            Xin%Xu%NULL%3,                            Jun%Huang%NULL%2,                            Jian-fu%Li%NULL%3,                            Li-min%Ou%NULL%3,                            Bo%Cheng%NULL%3,                            Shan%Xiong%NULL%3,                            Zhan-hong%Xie%NULL%2,                            Zheng-yi%Ni%NULL%0,                            Yu%Hu%NULL%0,                            Lei%Liu%NULL%0,                            Hong%Shan%NULL%0,                            Chun-liang%Lei%NULL%0,                            Yi-xiang%Peng%NULL%0,                            Li%Wei%NULL%0,                            Yong%Liu%NULL%0,                            Ya-hua%Hu%NULL%0,                            Peng%Peng%NULL%0,                            Jian-ming%Wang%NULL%0,                            Ji-yang%Liu%NULL%0,                            Zhong%Chen%NULL%0,                            Gang%Li%NULL%0,                            Zhi-jian%Zheng%NULL%0,                            Shao-qin%Qiu%NULL%0,                            Jie%Luo%NULL%0,                            Chang-jiang%Ye%NULL%0,                            Shao-yong%Zhu%NULL%0,                            Lin-ling%Cheng%NULL%3,                            Feng%Ye%NULL%5,                            Shi-yue%Li%NULL%0,                            Jin-ping%Zheng%NULL%3,                            Nuo-fu%Zhang%NULL%3,                            Nan-shan%Zhong%NULL%0,                            Jian-xing%He%NULL%0]"
$ws.Range("E26").Value = "[Christopher M%Petrilli%NULL%0,                            Simon A%Jones%NULL%1,                            Jie%Yang%NULL%0,                            Harish%Rajagopalan%NULL%1,                            Luke%O’Donnell%NULL%1,                            Yelena%Chernyak%NULL%1,                            Katie A%Tobin%NULL%1,                            Robert J%Cerfolio%NULL%1,                            Fritz%Francois%NULL%1,                            Leora I%Horwitz%NULL%1]"
$ws.Range("E27").Value = "[NULL%NULL%NULL%0,                            NULL%NULL%NULL%0,                            NULL%NULL%NULL%0,                            Stephanie%Bialek%NULL%1,                            Ellen%Boundy%NULL%1,                            Virginia%Bowen%NULL%1,                            Nancy%Chow%NULL%2,                            Amanda%Cohn%NULL%1,                            Nicole%Dowling%NULL%1,                            Sascha%Ellington%NULL%1,                            Ryan%Gierke%NULL%2,                            Aron%Hall%NULL%2,                            Jessica%MacNeil%NULL%1,                            Priti%Patel%NULL%1,                            Georgina%Peacock%NULL%1,                            Tamara%Pilishvili%NULL%2,                            Hilda%Razzaghi%NULL%1,                            Nia%Reed%NULL%1,                            Matthew%Ritchey%NULL%2,                            Erin%Sauber-Schatz%NULL%1]"
$ws.Range("E28").Value = "[Rodrigo B.%Serafim%NULL%1,                            Pedro%Póvoa%NULL%1,                            Vicente%Souza-Dantas%NULL%1,                            André C.%Kalil%NULL%1,                            Jorge I.F.%Salluh%NULL%1]"
$ws.Range("E29").Value = "[Ryan P%Barbaro%NULL%2,                            Graeme%MacLaren%NULL%1,                            Philip S%Boonstra%NULL%1,                            Theodore J%Iwashyna%NULL%1,                            Arthur S%Slutsky%NULL%1,                            Eddy%Fan%NULL%1,                            Robert H%Bartlett%NULL%1,                            Joseph E%Tonna%NULL%1,                            Robert%Hyslop%NULL%1,                            Jeffrey J%Fanning%NULL%1,                            Peter T%Rycus%NULL%1,                            Steve J%Hyer%NULL%1,                            Marc M%Anders%NULL%1,                            Cara L%Agerstrand%NULL%1,                            Katarzyna%Hryniewicz%NULL%1,                            Rodrigo%Diaz%NULL%1,                            Roberto%Lorusso%NULL%1,                            Alain%Combes%NULL%1,                            Daniel%Brodie%NULL%4,                            Peta%Alexander%NULL%1,                            Nicholas%Barrett%NULL%1,                            Jan%Bělohlávek%NULL%1,                            Dale%Fisher%NULL%1,                            John%Fraser%NULL%1,                            Ali Ait%Hssain%NULL%1,                            Jae Sung%Jung%NULL%1,                            Michael%McMullan%NULL%1,                            Yatin%Mehta%NULL%1,                            Mark T.%Ogino%NULL%1,                            Matthew L.%Paden%NULL%1,                            Kiran%Shekar%NULL%1,                            Christine%Stead%NULL%1,                            Yasir%Abu-Omar%NULL%1,                            Vanni%Agnoletti%NULL%1,                            Anzila%Akbar%NULL%1,                            Huda%Alfoudri%NULL%1,                            Carlos%Alviar%NULL%1,                            Vladimir%Aronsky%NULL%1,                            Erin%August%NULL%1,                            Georg%Auzinger%NULL%1,                            Hilda%Aveja%NULL%1,                            Rhonda%Bakken%NULL%1,                            Joan%Balcells%NULL%1,                            Sripal%Bangalore%NULL%1,                            Bernard W.%Barnes%NULL%1,                            Alaiza%Bautista%NULL%1,                            Lorraine L.%Bellows%NULL%1,                            Felipe%Beltran%NULL%1,                            Peyman%Benharash%NULL%1,                            Marco%Benni%NULL%1,                            Jennifer%Berg%NULL%1,                            Pietro%Bertini%NULL%1,                            Pablo%Blanco-Schweizer%NULL%1,                            Melissa%Brunsvold%NULL%1,                            Jenny%Budd%NULL%1,                            Debra%Camp%NULL%1,                            Mark%Caridi-Scheible%NULL%0,                            Edmund%Carton%NULL%1,                            Elena%Casanova-Ghosh%NULL%1,                            Anthony%Castleberry%NULL%1,                            Christopher T.%Chipongian%NULL%1,                            Chang Woo%Choi%NULL%1,                            Alessandro%Circelli%NULL%1,                            Elliott%Cohen%NULL%1,                            Michael%Collins%NULL%1,                            Scott%Copus%NULL%1,                            Jill%Coy%NULL%1,                            Brandon%Crist%NULL%1,                            Leonora%Cruz%NULL%1,                            Mirosław%Czuczwar%NULL%1,                            Mani%Daneshmand%NULL%1,                            Daniel%Davis II%NULL%1,                            Kim%De la Cruz%NULL%1,                            Cyndie%Devers%NULL%1,                            Toni%Duculan%NULL%1,                            Lucian%Durham%NULL%1,                            Subbarao%Elapavaluru%NULL%1,                            Carlos V.%Elzo Kraemer%NULL%1,                            EDMÍLSON CARDOSO%Filho%NULL%1,                            Jillian%Fitzgerald%NULL%1,                            Giuseppe%Foti%NULL%1,                            Matthew%Fox%NULL%1,                            David%Fritschen%NULL%1,                            David%Fullerton%NULL%1,                            Elton%Gelandt%NULL%1,                            Stacy%Gerle%NULL%1,                            Marco%Giani%NULL%1,                            Si Guim%Goh%NULL%1,                            Sara%Govener%NULL%1,                            Julie%Grone%NULL%1,                            Miles%Guber%NULL%1,                            Vadim%Gudzenko%NULL%1,                            Daniel%Gutteridge%NULL%1,                            Jennifer%Guy%NULL%1,                            Jonathan%Haft%NULL%1,                            Cameron%Hall%NULL%1,                            Ibrahim Fawzy%Hassan%NULL%1,                            Rubén%Herrán%NULL%1,                            Hitoshi%Hirose%NULL%1,                            Abdulsalam Saif%Ibrahim%NULL%1,                            Don%Igielski%NULL%1,                            Felicia A.%Ivascu%NULL%1,                            Jaume%Izquierdo Blasco%NULL%1,                            Julie%Jackson%NULL%1,                            Harsh%Jain%NULL%1,                            Bhavini%Jaiswal%NULL%1,                            Andrea C.%Johnson%NULL%1,                            Jenniver A.%Jurynec%NULL%1,                            Norma M%Kellter%NULL%1,                            Adam%Kohl%NULL%1,                            Zachary%Kon%NULL%1,                            Markus%Kredel%NULL%1,                            Karen%Kriska%NULL%1,                            Chandra%Kunavarapu%NULL%1,                            Oude%Lansink-Hartgring%NULL%1,                            Jeliene%LaRocque%NULL%1,                            Sharon Beth%Larson%NULL%1,                            Tracie%Layne%NULL%1,                            Stephane%Ledot%NULL%1,                            Napolitan%Lena%NULL%1,                            Jonathan%Lillie%NULL%1,                            Gösta%Lotz%NULL%1,                            Mark%Lucas%NULL%1,                            Lee%Ludwigson%NULL%1,                            Jacinta J.%Maas%NULL%1,                            Joanna%Maertens%NULL%1,                            David%Mast%NULL%1,                            Scott%McCardle%NULL%1,                            Bernard%McDonald%NULL%1,                            Allison%McLarty%NULL%1,                            Chelsea%McMahon%NULL%1,                            Patrick%Meybohm%NULL%1,                            Bart%Meyns%NULL%1,                            Casey%Miller%NULL%1,                            Fernando%Moraes Neto%NULL%1,                            Kelly%Morris%NULL%1,                            Ralf%Muellenbach%NULL%1,                            Meghan%Nicholson%NULL%1,                            Serena%O'Brien%NULL%1,                            Kathryn%O'Keefe%NULL%1,                            Tawnya%Ogston%NULL%1,                            Gary%Oldenburg%NULL%1,                            Fabiana M.%Oliveira%NULL%1,                            Emily%Oppel%NULL%1,                            Diego%Pardo%NULL%2,                            Diego%Pardo%NULL%0,                            Sara J.%Parker%NULL%1,                            Finn M.%Pedersen%NULL%1,                            Crescens%Pellecchia%NULL%1,                            Jose A.S.%Pelligrini%NULL%1,                            Thao T.N.%Pham%NULL%1,                            Ann R.%Phillips%NULL%1,                            Tasneem%Pirani%NULL%1,                            Paweł%Piwowarczyk%NULL%1,                            Robert%Plambeck%NULL%1,                            William%Pruett%NULL%1,                            Brittany%Quandt%NULL%1,                            Kollengode%Ramanathan%NULL%1,                            Alejandro%Rey%NULL%1,                            Christian%Reyher%NULL%1,                            Jordi%Riera del Brio%NULL%1,                            Rachel%Roberts%NULL%1,                            David%Roe%NULL%1,                            Peter P.%Roeleveld%NULL%1,                            Janet%Rudy%NULL%1,                            Luis F.%Rueda%NULL%1,                            Emanuele%Russo%NULL%1,                            Jesús%Sánchez Ballesteros%NULL%1,                            Nancy%Satou%NULL%1,                            Mauricio Guidi%Saueressig%NULL%1,                            Paul C.%Saunders%NULL%1,                            Margaret%Schlotterbeck%NULL%1,                            Patricia%Schwarz%NULL%1,                            Nicole%Scriven%NULL%1,                            Alexis%Serra%NULL%1,                            Mohammad%Shamsah%NULL%1,                            Lucy%Sim%NULL%1,                            Alexandra%Smart%NULL%1,                            Adam%Smith%NULL%1,                            Deane%Smith%NULL%1,                            Maggie%Smith%NULL%1,                            Neel%Sodha%NULL%1,                            Michael%Sonntagbauer%NULL%1,                            Marc%Sorenson%NULL%1,                            Eric B%Stallkamp%NULL%1,                            Allison%Stewart%NULL%1,                            Kathy%Swartz%NULL%1,                            Koji%Takeda%NULL%1,                            Shaun%Thompson%NULL%1,                            Bridget%Toy%NULL%1,                            Divina%Tuazon%NULL%1,                            Makoto%Uchiyama%NULL%1,                            Obiora I.%Udeozo%NULL%1,                            Scott%van Poppel%NULL%1,                            Corey%Ventetuolo%NULL%1,                            Leen%Vercaemst%NULL%1,                            Nguyen V.%Vinh Chau%NULL%1,                            I-Wen%Wang%NULL%1,                            Carrie%Williamson%NULL%1,                            Brock%Wilson%NULL%1,                            Helen%Winkels%NULL%1]"
